$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 2
$ws.Range("I2").Value = "18 TL - 18 TL"

# Row 3
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Row 14
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
